# insert_value_to_table.xlsx - "debug started, have done some fixes"
#
# Translate the profession names in column A (rows 2-11) from Russian to
# English, add a new "label" header column (AG1), mark the translated
# column as Text-formatted, and reset the active selection/scroll back to
# the top-left of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header column AG: "label" --------------------------------------
$ws.Range("AG1").Value = "label"

# --- Translate professions (column A, rows 2-11) -----------------------
$ws.Range("A2").Value  = "Farmer"
$ws.Range("A3").Value  = "Biologist"
$ws.Range("A4").Value  = "Assembler on the production line"
$ws.Range("A5").Value  = "Repair Specialist"
$ws.Range("A6").Value  = "Developer"
$ws.Range("A7").Value  = "Architect"
$ws.Range("A8").Value  = "Writer"
$ws.Range("A9").Value  = "Restorer"
$ws.Range("A10").Value = "Sales floor manager"
$ws.Range("A11").Value = "Advocate"

# Those cells now hold plain labels -> format them as Text (numFmtId 49)
$ws.Range("A2:A11").NumberFormat = "@"

# --- Reset selection / scroll position -----------------------------------
$ws.Range("A11").Select()
